$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -130.4
$ws.Range("B3").Value = -284.9
$ws.Range("C3").Value = -198.4
$ws.Range("C4").Value = -83.59999999999999
$ws.Range("C5").Value = 39.6
$ws.Range("C11").Value = 221.2
$ws.Range("C12").Value = 278.7
$ws.Range("C13").Value = 307
